$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update restated financial figures for 2014-2018 (rows 2-6) ---
# Row 2
$ws.Range("D2").Value = 8357
$ws.Range("E2").Value = 283
$ws.Range("F2").Value = 282
$ws.Range("G2").Value = 177
$ws.Range("H2").Value = 105
$ws.Range("I2").Value = 117
$ws.Range("J2").Value = -12
$ws.Range("K2").Value = 7275
$ws.Range("L2").Value = 4622
$ws.Range("M2").Value = 2652
$ws.Range("N2").Value = 2668
$ws.Range("O2").Value = -16
$ws.Range("P2").Value = 184
$ws.Range("Q2").Value = 639
$ws.Range("R2").Value = -526
$ws.Range("S2").Value = 1236
$ws.Range("T2").Value = 482
$ws.Range("U2").Value = 157
$ws.Range("V2").Value = 3656
$ws.Range("W2").Value = 3.39
$ws.Range("X2").Value = 1.26
$ws.Range("Y2").Value = 4.63
$ws.Range("Z2").Value = 1.63
$ws.Range("AA2").Value = 174.27
$ws.Range("AB2").Value = 1365.74
$ws.Range("AC2").Value = 318
$ws.Range("AD2").Value = 42.5
$ws.Range("AE2").Value = 7681
$ws.Range("AF2").Value = 1.76
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 36727943
# Row 3
$ws.Range("D3").Value = 8671
$ws.Range("E3").Value = 356
$ws.Range("F3").Value = 356
$ws.Range("G3").Value = 202
$ws.Range("H3").Value = 107
$ws.Range("I3").Value = 110
$ws.Range("J3").Value = -4
$ws.Range("K3").Value = 7761
$ws.Range("L3").Value = 5104
$ws.Range("M3").Value = 2658
$ws.Range("N3").Value = 2681
$ws.Range("O3").Value = -23
$ws.Range("P3").Value = 184
$ws.Range("Q3").Value = 413
$ws.Range("R3").Value = -1918
$ws.Range("S3").Value = 216
$ws.Range("T3").Value = 99
$ws.Range("U3").Value = 314
$ws.Range("V3").Value = 3940
$ws.Range("W3").Value = 4.11
$ws.Range("X3").Value = 1.23
$ws.Range("Y3").Value = 4.12
$ws.Range("Z3").Value = 1.42
$ws.Range("AA3").Value = 192.03
$ws.Range("AB3").Value = 1455.05
$ws.Range("AC3").Value = 300
$ws.Range("AD3").Value = 54.39
$ws.Range("AE3").Value = 7447
$ws.Range("AF3").Value = 2.19
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 36727943
# Row 4
$ws.Range("D4").Value = 9476
$ws.Range("E4").Value = 423
$ws.Range("F4").Value = 423
$ws.Range("G4").Value = 411
$ws.Range("H4").Value = 265
$ws.Range("I4").Value = 267
$ws.Range("J4").Value = -2
$ws.Range("K4").Value = 6631
$ws.Range("L4").Value = 3691
$ws.Range("M4").Value = 2940
$ws.Range("N4").Value = 2904
$ws.Range("O4").Value = 36
$ws.Range("P4").Value = 184
$ws.Range("Q4").Value = 332
$ws.Range("R4").Value = 1049
$ws.Range("S4").Value = -1592
$ws.Range("T4").Value = 301
$ws.Range("U4").Value = 31
$ws.Range("V4").Value = 2555
$ws.Range("W4").Value = 4.47
$ws.Range("X4").Value = 2.79
$ws.Range("Y4").Value = 9.56
$ws.Range("Z4").Value = 3.68
$ws.Range("AA4").Value = 125.52
$ws.Range("AB4").Value = 1597.13
$ws.Range("AC4").Value = 727
$ws.Range("AD4").Value = 16.24
$ws.Range("AE4").Value = 8068
$ws.Range("AF4").Value = 1.46
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 0.85
$ws.Range("AI4").Value = 13.49
$ws.Range("AJ4").Value = 36727943
# Row 5
$ws.Range("D5").Value = 10233
$ws.Range("E5").Value = 417
$ws.Range("F5").Value = 417
$ws.Range("G5").Value = 514
$ws.Range("H5").Value = 394
$ws.Range("I5").Value = 396
$ws.Range("J5").Value = -3
$ws.Range("K5").Value = 7253
$ws.Range("L5").Value = 4048
$ws.Range("M5").Value = 3205
$ws.Range("N5").Value = 3159
$ws.Range("O5").Value = 46
$ws.Range("P5").Value = 184
$ws.Range("Q5").Value = 416
$ws.Range("R5").Value = -380
$ws.Range("S5").Value = 152
$ws.Range("T5").Value = 547
$ws.Range("U5").Value = -131
$ws.Range("V5").Value = 2823
$ws.Range("W5").Value = 4.08
$ws.Range("X5").Value = 3.85
$ws.Range("Y5").Value = 13.07
$ws.Range("Z5").Value = 5.67
$ws.Range("AA5").Value = 126.3
$ws.Range("AB5").Value = 1796.6
$ws.Range("AC5").Value = 1079
$ws.Range("AD5").Value = 10.75
$ws.Range("AE5").Value = 8873
$ws.Range("AF5").Value = 1.31
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 0.86
$ws.Range("AI5").Value = 8.98
$ws.Range("AJ5").Value = 36727943
# Row 6
$ws.Range("D6").Value = 10157
$ws.Range("E6").Value = 237
$ws.Range("F6").Value = 237
$ws.Range("G6").Value = 71
$ws.Range("H6").Value = 35
$ws.Range("I6").Value = 37
$ws.Range("K6").Value = 8832
$ws.Range("L6").Value = 5733
$ws.Range("M6").Value = 3098
$ws.Range("N6").Value = 3058
$ws.Range("P6").Value = 184
$ws.Range("Q6").Value = -384
$ws.Range("R6").Value = -939
$ws.Range("S6").Value = 1513
$ws.Range("T6").Value = 483
$ws.Range("U6").Value = -867
$ws.Range("V6").Value = 4465
$ws.Range("W6").Value = 2.34
$ws.Range("X6").Value = 0.34
$ws.Range("Y6").Value = 1.2
$ws.Range("Z6").Value = 0.43
$ws.Range("AA6").Value = 185.05
$ws.Range("AB6").Value = 1787.28
$ws.Range("AC6").Value = 102
$ws.Range("AD6").Value = 76.5
$ws.Range("AE6").Value = 8713
$ws.Range("AF6").Value = 0.89
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 0.64
$ws.Range("AI6").Value = 46.98
$ws.Range("AJ6").Value = 36727943

# --- Clear forecast rows 7-9 (2019E-2021E) data beyond column C: data no longer available ---
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
